$d = $word.ActiveDocument

# 1. Title heading and the later bold "Play ..." repeat (both identical text)
$d.Content.Find.Execute(
    "Play Inferno Gladiator Free: Review & Pros and Cons 2021", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Inferno Gladiator Free - Exciting Gameplay Features", 2)

# 2. "What we like" bullet 1
$d.Content.Find.Execute(
    "Well-crafted graphics and symbols", $true, $false, $false, $false, $false,
    $true, 1, $false, "Exciting gameplay features with Wild and Scatter symbols", 2)

# 3. "What we like" bullet 2
$d.Content.Find.Execute(
    "Free Spins round with multiplier feature", $true, $false, $false, $false, $false,
    $true, 1, $false, "Free spins round with a fixed Wild and increasing multiplier", 2)

# 4. "What we like" bullet 3
$d.Content.Find.Execute(
    "Decent potential winnings up to 2,000x the bet placed", $true, $false, $false, $false, $false,
    $true, 1, $false, "Well-crafted graphics inspired by ancient Rome", 2)

# 5. "What we like" bullet 4
$d.Content.Find.Execute(
    "Available on a range of different online casinos", $true, $false, $false, $false, $false,
    $true, 1, $false, "Decent win potential with a maximum jackpot of 250,000", 2)

# 6. "What we don't like" bullet 1
$d.Content.Find.Execute(
    "Low RTP at 95.04%", $true, $false, $false, $false, $false,
    $true, 1, $false, "Some players may find the game falls short of expectations", 2)

# 7. "What we don't like" bullet 2
$d.Content.Find.Execute(
    "May not meet high expectations", $true, $false, $false, $false, $false,
    $true, 1, $false, "Low RTP of 95.04%", 2)

# 8. Italic summary paragraph
$d.Content.Find.Execute(
    "Read our review of Inferno Gladiator and learn about its pros and cons, gameplay features, win potential, and availability to play for free in 2021 on a range of online casinos.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Read our review of Inferno Gladiator and play it for free. Exciting gameplay features and decent win potential.", 2)
